$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 3.42 = 13157.11 pesos
✅ 13157.11 pesos = 3.39 = 916.44 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%"

# --- tasas: update the N/O rate columns ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 292.777
$wsTasas.Range("O10").Value = 3852.1
$wsTasas.Range("N12").Value = 3883.5
$wsTasas.Range("O12").Value = 270.5
